# ECP-962: adds property frequency on turnover
#
# Inserts a new "frequency" column (with MONTHLY / QUARTERLY sample
# values) between the existing "type" and "turnover gross amount"
# columns on the TurnoverImport fixture sheet. Everything from the old
# "turnover gross amount" column onward shifts one column to the right
# (E->F, F->G, G->H, H->I, I->J, J->K), while the column widths/styles
# defined on the sheet (<cols>) are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot the current (pre-edit) values for the columns that are
#     going to be shifted, so later writes don't clobber data we still
#     need to read. ---
$srcCols = @("E", "F", "G", "H", "I", "J")
$rows = @(1, 2, 3)

# NB: read with Value2 (not Value) - in this runtime plain `.Value`
# returns the property descriptor instead of the cell's content, and
# `.Text` would give back a locale-formatted display string (losing
# numeric precision). `.Value2` gives the raw underlying value and
# round-trips cleanly through `.Value2 = ...` on write.
$values = @{}
foreach ($col in $srcCols) {
    foreach ($row in $rows) {
        $addr = "$col$row"
        $values[$addr] = $ws.Range($addr).Value2
    }
}

# Only columns E and F originally carry the "currency-ish" style
# (style index 2 - numFmtId 44). Capture its NumberFormat once; every
# other column is left at the default "General" style. (Re-assigning
# NumberFormat = "General" to a default cell would actually *create* a
# brand-new, unwanted number format, so we simply never touch
# NumberFormat for those cells.)
$style2Format = $ws.Range("E2").NumberFormat
$styledSrcCols = @{ "E" = $true; "F" = $true }

# --- shift E..J one column to the right, into F..K, right-to-left so
#     we never overwrite a source cell before it has been read. ---
$destCols = @("F", "G", "H", "I", "J", "K")
for ($i = $srcCols.Length - 1; $i -ge 0; $i--) {
    $srcCol = $srcCols[$i]
    $dstCol = $destCols[$i]
    foreach ($row in $rows) {
        $srcAddr = "$srcCol$row"
        $dstAddr = "$dstCol$row"
        if ($styledSrcCols.ContainsKey($srcCol)) {
            $ws.Range($dstAddr).NumberFormat = $style2Format
        }
        $ws.Range($dstAddr).Value2 = $values[$srcAddr]
    }
}

# Row 3 had no value in the old "turnover net amount" column (old F3),
# which after the shift lands on G3: make sure it stays blank but keeps
# the currency-style formatting (matches source: <c r="G3" s="2"/>).
$ws.Range("G3").NumberFormat = $style2Format
$ws.Range("G3").Value2 = ""

# --- populate the new "frequency" column (E) ---
$ws.Range("E1").NumberFormat = $style2Format
$ws.Range("E1").Value2 = "frequency"

$ws.Range("E2").NumberFormat = $style2Format
$ws.Range("E2").Value2 = "MONTHLY"

$ws.Range("E3").NumberFormat = $style2Format
$ws.Range("E3").Value2 = "QUARTERLY"

# Match the saved selection in the authored workbook.
$ws.Range("E3").Select()
